$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 63287.8
$ws.Range("I21").Value = 58012.332
$ws.Range("K21").Value = 58012.332
$ws.Range("M21").Value = -57544.332
$ws.Range("H23").Value = 63287.8
$ws.Range("I23").Value = 58012.332
$ws.Range("K23").Value = 58012.332
$ws.Range("M23").Value = -57778.332
$ws.Range("H29").Value = 2422.9
$ws.Range("I29").Value = 1886.5
$ws.Range("J29").Value = 3227.5
$ws.Range("K29").Value = 5659.5
$ws.Range("L29").Value = 9682.5
$ws.Range("M29").Value = -5378.5
$ws.Range("N29").Value = -10244.5
$ws.Range("H38").Value = 834016.8
$ws.Range("I38").Value = 2500056.8
$ws.Range("J38").Value = 996.875
$ws.Range("K38").Value = 7500170.399999999
$ws.Range("L38").Value = 2990.625
$ws.Range("M38").Value = -7499798.399999999
$ws.Range("N38").Value = -3734.625
$ws.Range("H58").Value = 3312.3
$ws.Range("I58").Value = 185
$ws.Range("J58").Value = 4652.5713
$ws.Range("K58").Value = 555
$ws.Range("L58").Value = 13957.7139
$ws.Range("M58").Value = -405
$ws.Range("N58").Value = -14257.7139
$ws.Range("H86").Value = 14288687
$ws.Range("I86").Value = 20001760
$ws.Range("J86").Value = 6002
$ws.Range("K86").Value = 20001760
$ws.Range("L86").Value = 6002
$ws.Range("M86").Value = -20000637
$ws.Range("N86").Value = -8248
$ws.Range("H87").Value = 42177
$ws.Range("J87").Value = 42177
$ws.Range("L87").Value = 42177
$ws.Range("N87").Value = -44673
$ws.Range("H88").Value = 3973.5
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 894
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 894
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -1706
$ws.Range("H89").Value = 14288687
$ws.Range("I89").Value = 20001760
$ws.Range("J89").Value = 6002
$ws.Range("K89").Value = 100008800
$ws.Range("L89").Value = 30010
$ws.Range("M89").Value = -100003184
$ws.Range("N89").Value = -41242
$ws.Range("H90").Value = 42177
$ws.Range("J90").Value = 42177
$ws.Range("L90").Value = 126531
$ws.Range("N90").Value = -139011
$ws.Range("H91").Value = 3973.5
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 894
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 894
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -3702
$ws.Range("H138").Value = 6343678.5
$ws.Range("I138").Value = 1414.7858
$ws.Range("J138").Value = 9825706
$ws.Range("K138").Value = 4244.357400000001
$ws.Range("L138").Value = 29477118
$ws.Range("M138").Value = 895.6425999999992
$ws.Range("N138").Value = -29487398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2039259
$ws.Range("I132").Value = 2041040.6
$ws.Range("K132").Value = 6123121.800000001
$ws.Range("M132").Value = -6120591.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 390578.38
$ws.Range("I86").Value = 1818.5454
$ws.Range("J86").Value = 1001486.7
$ws.Range("K86").Value = 1818.5454
$ws.Range("L86").Value = 1001486.7
$ws.Range("M86").Value = -695.5454
$ws.Range("N86").Value = -1003732.7
$ws.Range("H89").Value = 390578.38
$ws.Range("I89").Value = 1818.5454
$ws.Range("J89").Value = 1001486.7
$ws.Range("K89").Value = 9092.726999999999
$ws.Range("L89").Value = 5007433.5
$ws.Range("M89").Value = -3476.726999999999
$ws.Range("N89").Value = -5018665.5
$ws.Range("H132").Value = 1000000000
$ws.Range("J132").Value = 1000000000
$ws.Range("L132").Value = 1000000000
$ws.Range("N132").Value = -1000010120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 820.9231
$ws.Range("I105").Value = 831
$ws.Range("K105").Value = 831
$ws.Range("M105").Value = 916

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 503.3793
$ws.Range("I122").Value = 232.61111
$ws.Range("J122").Value = 946.4545000000001
$ws.Range("K122").Value = 2093.49999
$ws.Range("L122").Value = 8518.0905
$ws.Range("M122").Value = 356.5000100000002
$ws.Range("N122").Value = -13418.0905
$ws.Range("H132").Value = 1762.5
$ws.Range("I132").Value = 866.6667
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 7800.0003
$ws.Range("L132").Value = 20700
$ws.Range("M132").Value = -5270.0003
$ws.Range("N132").Value = -25760
$ws.Range("H138").Value = 8549388
$ws.Range("I138").Value = 1053.8889
$ws.Range("J138").Value = 15876531
$ws.Range("K138").Value = 3161.6667
$ws.Range("L138").Value = 47629593
$ws.Range("M138").Value = 1978.3333
$ws.Range("N138").Value = -47639873
$ws.Range("H141").Value = 4390712.5
$ws.Range("I141").Value = 1414.4445
$ws.Range("J141").Value = 5752908.5
$ws.Range("K141").Value = 4243.333500000001
$ws.Range("L141").Value = 17258725.5
$ws.Range("M141").Value = 936.6664999999994
$ws.Range("N141").Value = -17269085.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3337896.8
$ws.Range("I70").Value = 5267394.5
$ws.Range("J70").Value = 5127.727
$ws.Range("K70").Value = 5267394.5
$ws.Range("L70").Value = 5127.727
$ws.Range("M70").Value = -5267124.5
$ws.Range("N70").Value = -5667.727
$ws.Range("H73").Value = 3337896.8
$ws.Range("I73").Value = 5267394.5
$ws.Range("J73").Value = 5127.727
$ws.Range("K73").Value = 5267394.5
$ws.Range("L73").Value = 5127.727
$ws.Range("M73").Value = -5266458.5
$ws.Range("N73").Value = -6999.727
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 754.8182
$ws.Range("I16").Value = 567
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 567
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -397
$ws.Range("N16").Value = -1940
$ws.Range("H82").Value = 1930.5264
$ws.Range("I82").Value = 1399.091
$ws.Range("J82").Value = 2661.25
$ws.Range("K82").Value = 1399.091
$ws.Range("L82").Value = 2661.25
$ws.Range("M82").Value = -1038.091
$ws.Range("N82").Value = -3383.25
$ws.Range("H85").Value = 1930.5264
$ws.Range("I85").Value = 1399.091
$ws.Range("J85").Value = 2661.25
$ws.Range("K85").Value = 1399.091
$ws.Range("L85").Value = 2661.25
$ws.Range("M85").Value = -151.0909999999999
$ws.Range("N85").Value = -5157.25
$ws.Range("H122").Value = 2514.9697
$ws.Range("I122").Value = 2250.1904
$ws.Range("J122").Value = 2978.3333
$ws.Range("K122").Value = 6750.5712
$ws.Range("L122").Value = 8934.999899999999
$ws.Range("M122").Value = -4300.5712
$ws.Range("N122").Value = -13834.9999
$ws.Range("H132").Value = 242142.27
$ws.Range("I132").Value = 55875.434
$ws.Range("J132").Value = 772286.3
$ws.Range("K132").Value = 167626.302
$ws.Range("L132").Value = 2316858.9
$ws.Range("M132").Value = -165096.302
$ws.Range("N132").Value = -2321918.9
$ws.Range("H134").Value = 36762.438
$ws.Range("J134").Value = 36762.438
$ws.Range("L134").Value = 36762.438
$ws.Range("N134").Value = -46902.438
$ws.Range("H135").Value = 49972.727
$ws.Range("J135").Value = 49972.727
$ws.Range("L135").Value = 49972.727
$ws.Range("N135").Value = -60112.727
$ws.Range("H141").Value = 53000
$ws.Range("J141").Value = 53000
$ws.Range("L141").Value = 53000
$ws.Range("N141").Value = -63360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2380.9539
$ws.Range("I132").Value = 556.58
$ws.Range("J132").Value = 8462.200000000001
$ws.Range("K132").Value = 1669.74
$ws.Range("L132").Value = 25386.6
$ws.Range("M132").Value = 860.2599999999998
$ws.Range("N132").Value = -30446.6
